$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New table rows (8-10) describing additional power-supply / DCLink devices
$ws.Range("A8").Value = "10.128.102.104"
$ws.Range("B8").Value = "BBB-PS-DCLink"
$ws.Range("C8").Value = "PowerSupply"
$ws.Range("D8").Value = "1,2 "
$ws.Range("E8").Value = "FBP-DCLink"

$ws.Range("A9").Value = "10.128.102.122"
$ws.Range("B9").Value = "BBB-SI-CORRETORAS1"
$ws.Range("C9").Value = "PowerSupply"
$ws.Range("D9").Value = "1,2,3,4,5,6,7"
$ws.Range("E9").Value = "CH1,CH2,CV1,CV2,CH3,CH4,CV3"

$ws.Range("A10").Value = "10.128.102.132"
$ws.Range("B10").Value = "BBB-SI-CORRETORAS2"
$ws.Range("C10").Value = "PowerSupply"
$ws.Range("D10").Value = "1,2,3,4,5,6,7"
$ws.Range("E10").Value = "CV4,CH5,CH6,CV5,CV6,CH7,CV7"

# Column E needs to be wider to fit the new device-name values
$ws.Columns("E").ColumnWidth = 29.558559

# Move the active selection up one row
[void]$ws.Range("B12").Select()
